# Auto-generated from upstream diff: refresh Market Board price snapshots
# (currentAveragePrice* / Leve profit columns) across the 8 job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 667689.25
$ws.Range("I6").Value = 1001305
$ws.Range("K6").Value = 3003915
$ws.Range("M6").Value = -3003803
$ws.Range("H9").Value = 363.91666
$ws.Range("J9").Value = 341
$ws.Range("L9").Value = 341
$ws.Range("N9").Value = -679
$ws.Range("H19").Value = 1475.625
$ws.Range("I19").Value = 1474.6923
$ws.Range("K19").Value = 1474.6923
$ws.Range("M19").Value = -1299.6923
$ws.Range("H32").Value = 2512.375
$ws.Range("J32").Value = 2728.4285
$ws.Range("L32").Value = 2728.4285
$ws.Range("N32").Value = -3380.4285
$ws.Range("H70").Value = 2279.8
$ws.Range("I70").Value = 2150
$ws.Range("J70").Value = 2474.5
$ws.Range("K70").Value = 6450
$ws.Range("L70").Value = 7423.5
$ws.Range("M70").Value = -6180
$ws.Range("N70").Value = -7963.5
$ws.Range("H73").Value = 2279.8
$ws.Range("I73").Value = 2150
$ws.Range("J73").Value = 2474.5
$ws.Range("K73").Value = 6450
$ws.Range("L73").Value = 7423.5
$ws.Range("M73").Value = -5514
$ws.Range("N73").Value = -9295.5
$ws.Range("H103").Value = 850.2727
$ws.Range("J103").Value = 971
$ws.Range("L103").Value = 2913
$ws.Range("N103").Value = -4085
$ws.Range("H111").Value = 5437.5
$ws.Range("I111").Value = 876
$ws.Range("K111").Value = 2628
$ws.Range("M111").Value = 439
$ws.Range("H116").Value = 12705.866
$ws.Range("I116").Value = 20699.125
$ws.Range("J116").Value = 3570.7144
$ws.Range("K116").Value = 20699.125
$ws.Range("L116").Value = 3570.7144
$ws.Range("M116").Value = -17257.125
$ws.Range("N116").Value = -10454.7144
$ws.Range("H132").Value = 1659.093
$ws.Range("I132").Value = 1377.3235
$ws.Range("J132").Value = 2723.5557
$ws.Range("K132").Value = 4131.970499999999
$ws.Range("L132").Value = 8170.6671
$ws.Range("M132").Value = -1601.970499999999
$ws.Range("N132").Value = -13230.6671
$ws.Range("H137").Value = 2090066.1
$ws.Range("J137").Value = 4276429
$ws.Range("L137").Value = 12829287
$ws.Range("N137").Value = -12834387
$ws.Range("H138").Value = 2521.4
$ws.Range("J138").Value = 2737.7908
$ws.Range("L138").Value = 8213.3724
$ws.Range("N138").Value = -18493.3724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14706876
$ws.Range("I32").Value = 15152403
$ws.Range("K32").Value = 15152403
$ws.Range("M32").Value = -15152116
$ws.Range("H63").Value = 3626.5833
$ws.Range("J63").Value = 4753.4287
$ws.Range("L63").Value = 4753.4287
$ws.Range("N63").Value = -6125.4287
$ws.Range("H66").Value = 3626.5833
$ws.Range("J66").Value = 4753.4287
$ws.Range("L66").Value = 23767.1435
$ws.Range("N66").Value = -30631.1435
$ws.Range("H74").Value = 1751.7435
$ws.Range("I74").Value = 1636.5483
$ws.Range("K74").Value = 1636.5483
$ws.Range("M74").Value = -762.5482999999999
$ws.Range("H77").Value = 1751.7435
$ws.Range("I77").Value = 1636.5483
$ws.Range("K77").Value = 8182.7415
$ws.Range("M77").Value = -3814.7415

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H76").Value = 42294.332
$ws.Range("J76").Value = 42294.332
$ws.Range("L76").Value = 42294.332
$ws.Range("N76").Value = -42924.332
$ws.Range("H79").Value = 42294.332
$ws.Range("J79").Value = 42294.332
$ws.Range("L79").Value = 42294.332
$ws.Range("N79").Value = -44478.332
$ws.Range("H86").Value = 730.05
$ws.Range("I86").Value = 725.2308
$ws.Range("J86").Value = 739
$ws.Range("K86").Value = 725.2308
$ws.Range("L86").Value = 739
$ws.Range("M86").Value = 397.7692
$ws.Range("N86").Value = -2985
$ws.Range("H89").Value = 730.05
$ws.Range("I89").Value = 725.2308
$ws.Range("J89").Value = 739
$ws.Range("K89").Value = 3626.154
$ws.Range("L89").Value = 3695
$ws.Range("M89").Value = 1989.846
$ws.Range("N89").Value = -14927

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 504.85715
$ws.Range("I22").Value = 504.85715
$ws.Range("K22").Value = 504.85715
$ws.Range("M22").Value = -154.85715
$ws.Range("H141").Value = 641990.25
$ws.Range("J141").Value = 641990.25
$ws.Range("L141").Value = 641990.25
$ws.Range("N141").Value = -652350.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 2094.2307
$ws.Range("I114").Value = 639.2
$ws.Range("J114").Value = 3003.625
$ws.Range("K114").Value = 1917.6
$ws.Range("L114").Value = 9010.875
$ws.Range("M114").Value = 1336.4
$ws.Range("N114").Value = -15518.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12608.333
$ws.Range("I70").Value = 63000
$ws.Range("J70").Value = 4480.645
$ws.Range("K70").Value = 63000
$ws.Range("L70").Value = 4480.645
$ws.Range("M70").Value = -62730
$ws.Range("N70").Value = -5020.645
$ws.Range("H73").Value = 12608.333
$ws.Range("I73").Value = 63000
$ws.Range("J73").Value = 4480.645
$ws.Range("K73").Value = 63000
$ws.Range("L73").Value = 4480.645
$ws.Range("M73").Value = -62064
$ws.Range("N73").Value = -6352.645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3955.4443
$ws.Range("I22").Value = 3849
$ws.Range("J22").Value = 3968.75
$ws.Range("K22").Value = 3849
$ws.Range("L22").Value = 3968.75
$ws.Range("M22").Value = -3554
$ws.Range("N22").Value = -4558.75
$ws.Range("H27").Value = 3955.4443
$ws.Range("I27").Value = 3849
$ws.Range("J27").Value = 3968.75
$ws.Range("K27").Value = 3849
$ws.Range("L27").Value = 3968.75
$ws.Range("M27").Value = -3742
$ws.Range("N27").Value = -4182.75
$ws.Range("H40").Value = 2744.75
$ws.Range("I40").Value = 2666.6667
$ws.Range("K40").Value = 2666.6667
$ws.Range("M40").Value = -2530.6667
$ws.Range("H46").Value = 7245.4287
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 7443.407
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 7443.407
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -7819.407

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6028.5713
$ws.Range("I62").Value = 3250
$ws.Range("J62").Value = 7140
$ws.Range("K62").Value = 3250
$ws.Range("L62").Value = 7140
$ws.Range("M62").Value = -2626
$ws.Range("N62").Value = -8388
$ws.Range("H65").Value = 6028.5713
$ws.Range("I65").Value = 3250
$ws.Range("J65").Value = 7140
$ws.Range("K65").Value = 16250
$ws.Range("L65").Value = 35700
$ws.Range("M65").Value = -13130
$ws.Range("N65").Value = -41940
